$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A12 held a stray "10.6" number; repair it as a real date using the same
# date formatting already used by A3:A11 (copy format, then set the value).
$ws.Range("A3").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A12").Value = 41435
$ws.Range("A13").Value = 41436

# New run result row, entered in the same field order the workbook
# author used (filename, test/status note, then result file).
$ws.Range("B13").Value = "11062013_20d_Cvap_Dillon"
$ws.Range("D13").Value = "ajo aloitettu 16:18"
$ws.Range("C13").Value = "run_20130612T042435"

$ws.Range("C16").Select()
